$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Tipo", matching the style of the other header cells (A1:C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "Tipo"

# Update existing MSE (col B) and R2 (col C) values for rows 2-4
$ws.Range("B2").Value = 0.3421150429834521
$ws.Range("C2").Value = 0.9976345621564433

$ws.Range("B3").Value = 0.3421150429834521
$ws.Range("C3").Value = 0.9976345621564433

$ws.Range("B4").Value = 0.3421150429834521
$ws.Range("C4").Value = 0.9976345621564433

# Add new "Tipo" column data for rows 2-4
$ws.Range("D2").Value = "single"
$ws.Range("D3").Value = "single"
$ws.Range("D4").Value = "single"
